$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Michigan State_GAME_SCORES")

$ws.Range("O4:O27").Formula = '=AVERAGEIF(B4:M4, "<>0")'

$ws.Range("O6").Select() | Out-Null
